$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (loading_percent values) for rows 2-25
$colB = New-Object "object[,]" 24,1
$colB[0,0] = 5.896476547055094
$colB[1,0] = 5.785638564604787
$colB[2,0] = 5.717245135427095
$colB[3,0] = 5.689325584832277
$colB[4,0] = 5.684687706963098
$colB[5,0] = 5.716868750872302
$colB[6,0] = 5.858351426795625
$colB[7,0] = 6.131564099018144
$colB[8,0] = 6.327726517267261
$colB[9,0] = 6.415583513738482
$colB[10,0] = 6.44862416940033
$colB[11,0] = 6.441518921937131
$colB[12,0] = 6.418306544827898
$colB[13,0] = 6.404057606563406
$colB[14,0] = 6.321954457026194
$colB[15,0] = 6.271211419949322
$colB[16,0] = 6.241897360532192
$colB[17,0] = 6.231951103047054
$colB[18,0] = 6.27662657909217
$colB[19,0] = 6.425131030025209
$colB[20,0] = 6.520840228364817
$colB[21,0] = 6.469891349962448
$colB[22,0] = 6.274178823148308
$colB[23,0] = 6.058316741456117
$ws.Range("B2:B25").Value = $colB

# Columns D:I for rows 2-25
$blockDI = New-Object "object[,]" 24,6
$blockDI[0,0] = 4.322891434097585
$blockDI[0,1] = 16.50001308497689
$blockDI[0,2] = 25.14293919829933
$blockDI[0,3] = 31.01252376501491
$blockDI[0,4] = 14.25255102788124
$blockDI[0,5] = 19.30882781077258
$blockDI[1,0] = 4.333149428401794
$blockDI[1,1] = 15.56125721325922
$blockDI[1,2] = 24.98891542917677
$blockDI[1,3] = 30.64890075880527
$blockDI[1,4] = 14.27466603116684
$blockDI[1,5] = 19.43182168715426
$blockDI[2,0] = 4.339743961878594
$blockDI[2,1] = 14.9603928210345
$blockDI[2,2] = 24.90399771483853
$blockDI[2,3] = 30.43891790759098
$blockDI[2,4] = 14.29224360263695
$blockDI[2,5] = 19.51121530729685
$blockDI[3,0] = 4.342506160678126
$blockDI[3,1] = 14.70965405265451
$blockDI[3,2] = 24.87184543258717
$blockDI[3,3] = 30.35678982236842
$blockDI[3,4] = 14.30040715734994
$blockDI[3,5] = 19.54454490333466
$blockDI[4,0] = 4.342969356395846
$blockDI[4,1] = 14.66767277918694
$blockDI[4,2] = 24.86665531551374
$blockDI[4,3] = 30.34336314585488
$blockDI[4,4] = 14.30182298623162
$blockDI[4,5] = 19.55013825847644
$blockDI[5,0] = 4.339780910180802
$blockDI[5,1] = 14.9570346886598
$blockDI[5,2] = 24.9035541394155
$blockDI[5,3] = 30.43779623824654
$blockDI[5,4] = 14.29234965500886
$blockDI[5,5] = 19.51166084770138
$blockDI[6,0] = 4.326367187223335
$blockDI[6,1] = 16.18155400818645
$blockDI[6,2] = 25.08784472846433
$blockDI[6,3] = 30.88445359994888
$blockDI[6,4] = 14.25934342976517
$blockDI[6,5] = 19.35043287811136
$blockDI[7,0] = 4.30239388995656
$blockDI[7,1] = 18.46864013910899
$blockDI[7,2] = 25.52449493645684
$blockDI[7,3] = 31.86080577351706
$blockDI[7,4] = 14.22655352699934
$blockDI[7,5] = 19.06493310733132
$blockDI[8,0] = 4.286177331006751
$blockDI[8,1] = 20.11932269656972
$blockDI[8,2] = 25.88905656527285
$blockDI[8,3] = 32.63219580359812
$blockDI[8,4] = 14.22218667481072
$blockDI[8,5] = 18.87376245808892
$blockDI[9,0] = 4.279098380865274
$blockDI[9,1] = 20.82837373078346
$blockDI[9,2] = 26.06386338140962
$blockDI[9,3] = 32.99321292603224
$blockDI[9,4] = 14.22452795080088
$blockDI[9,5] = 18.79080715016488
$blockDI[10,0] = 4.276460271147108
$blockDI[10,1] = 21.09089311486526
$blockDI[10,2] = 26.13129632593397
$blockDI[10,3] = 33.13123149190105
$blockDI[10,4] = 14.22603969749142
$blockDI[10,5] = 18.75996913709966
$blockDI[11,0] = 4.277026548460197
$blockDI[11,1] = 21.03462039395894
$blockDI[11,2] = 26.11671926067938
$blockDI[11,3] = 33.10145080838195
$blockDI[11,4] = 14.22568627236673
$blockDI[11,5] = 18.76658508326643
$blockDI[12,0] = 4.278880491315808
$blockDI[12,1] = 20.85009106101861
$blockDI[12,2] = 26.06938663946466
$blockDI[12,3] = 33.00454238936344
$blockDI[12,4] = 14.22463977973413
$blockDI[12,5] = 18.78825856354478
$blockDI[13,0] = 4.280021614873063
$blockDI[13,1] = 20.73628342317213
$blockDI[13,2] = 26.04055359382106
$blockDI[13,3] = 32.945349375663
$blockDI[13,4] = 14.22408026148387
$blockDI[13,5] = 18.80160909196525
$blockDI[14,0] = 4.286645926402411
$blockDI[14,1] = 20.07214540044766
$blockDI[14,2] = 25.87780858385299
$blockDI[14,3] = 32.60879358687357
$blockDI[14,4] = 14.22212100652048
$blockDI[14,5] = 18.87926435741556
$blockDI[15,0] = 4.290785827383699
$blockDI[15,1] = 19.65402407819721
$blockDI[15,2] = 25.78023029944798
$blockDI[15,3] = 32.40481431844704
$blockDI[15,4] = 14.22202969236724
$blockDI[15,5] = 18.92792935332456
$blockDI[16,0] = 4.293195060444285
$blockDI[16,1] = 19.40959026216542
$blockDI[16,2] = 25.72495309232769
$blockDI[16,3] = 32.28845104028186
$blockDI[16,4] = 14.22238438586011
$blockDI[16,5] = 18.95629759986746
$blockDI[17,0] = 4.294015616795392
$blockDI[17,1] = 19.32615113234434
$blockDI[17,2] = 25.70638422170957
$blockDI[17,3] = 32.24922171669128
$blockDI[17,4] = 14.22257432376046
$blockDI[17,5] = 18.96596746047731
$blockDI[18,0] = 4.290342225037677
$blockDI[18,1] = 19.698941477071
$blockDI[18,2] = 25.7905303541514
$blockDI[18,3] = 32.42642992302969
$blockDI[18,4] = 14.22199724577401
$blockDI[18,5] = 18.92270982778834
$blockDI[19,0] = 4.278334791453125
$blockDI[19,1] = 20.90445388998346
$blockDI[19,2] = 26.08325622419185
$blockDI[19,3] = 33.03297236627063
$blockDI[19,4] = 14.22493017350184
$blockDI[19,5] = 18.78187693072913
$blockDI[20,0] = 4.270735029210885
$blockDI[20,1] = 21.65746963708246
$blockDI[20,2] = 26.28175581193567
$blockDI[20,3] = 33.43694641464811
$blockDI[20,4] = 14.23049183776208
$blockDI[20,5] = 18.6931876509171
$blockDI[21,0] = 4.274768594680775
$blockDI[21,1] = 21.25874862903668
$blockDI[21,2] = 26.17517323063359
$blockDI[21,3] = 33.22069376705403
$blockDI[21,4] = 14.22718917362664
$blockDI[21,5] = 18.7402163383622
$blockDI[22,0] = 4.290542686800877
$blockDI[22,1] = 19.67864694001652
$blockDI[22,2] = 25.78587113547612
$blockDI[22,3] = 32.41665467159402
$blockDI[22,4] = 14.22201064668288
$blockDI[22,5] = 18.92506836006318
$blockDI[23,0] = 4.308632342604367
$blockDI[23,1] = 17.82341123475683
$blockDI[23,2] = 25.39851427554646
$blockDI[23,3] = 31.58663725621944
$blockDI[23,4] = 14.2319769844385
$blockDI[23,5] = 19.13889673055544
$ws.Range("D2:I25").Value = $blockDI

# Column K for rows 2-25
$colK = New-Object "object[,]" 24,1
$colK[0,0] = 11.95290194740078
$colK[1,0] = 11.41519057254376
$colK[2,0] = 11.07551643021965
$colK[3,0] = 10.93450628504019
$colK[4,0] = 10.9109416005214
$colK[5,0] = 11.07362492131407
$colK[6,0] = 11.76822923031436
$colK[7,0] = 13.27519436624748
$colK[8,0] = 14.29845547617277
$colK[9,0] = 14.73873763053355
$colK[10,0] = 14.90182725581117
$colK[11,0] = 14.86686481022838
$colK[12,0] = 14.75222802295567
$colK[13,0] = 14.68153596161366
$colK[14,0] = 14.26917373722734
$colK[15,0] = 14.0097376524422
$colK[16,0] = 13.8581448166723
$colK[17,0] = 13.80641111513247
$colK[18,0] = 14.03760066410544
$colK[19,0] = 14.78599836195596
$colK[20,0] = 15.25392623088918
$colK[21,0] = 15.00612595888344
$colK[22,0] = 14.02501139446602
$colK[23,0] = 12.87612731785998
$ws.Range("K2:K25").Value = $colK

